# Apply the "Add files via upload" revision to WS_holdings.xlsx:
#   1) Bump the "as of" date in the confidential disclosure banner
#      from 2021-03-29 to 2021-03-30.
#   2) Refresh the Weight (D) / Percent Change (E) figures for rows 2-13
#      with the latest values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (no interactive password), which blocks direct
# cell writes. Unprotect for the duration of the edit, then restore
# protection once all the values are in place.
$ws.Unprotect()

# --- 1) Update the "as of" date in the disclosure text (cell A16) -----
$disclosure = $ws.Range("A16").Value()
$ws.Range("A16").Value = $disclosure.Replace("2021-03-29", "2021-03-30")

# --- 2) Update Weight / Percent Change figures, rows 2-13 -------------
$values = @{
    2  = @{ D = 0.03064983773323513;  E = 0 }
    3  = @{ D = 0.02436374617732993;  E = -0.007694157249338862 }
    4  = @{ D = 0.05354831658433233;  E = -0.009021842355175558 }
    5  = @{ D = 0.1386091040394344;   E = -0.006118286879673751 }
    6  = @{ D = 0.03130220326449595;  E = -0.006701414743112366 }
    7  = @{ D = 0.1186245827266151;   E = 0.01120896717373898 }
    8  = @{ D = 0.1019315533796987;   E = 0.005989180834621299 }
    9  = @{ D = 0.02819168851513697;  E = 0.0006764374295378239 }
    10 = @{ D = 0.121922805250158;    E = 0.007769372316499767 }
    11 = @{ D = 0.2483725337130391;   E = -0.007151029748283744 }
    12 = @{ D = 0.1024836286165243;   E = 0.000557795135944561 }
    13 = @{ D = 0.9999999999999998;   E = -0.0005408576407591736 }
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 4).Value = $values[$row].D
    $ws.Cells.Item($row, 5).Value = $values[$row].E
}

# Restore the worksheet protection that was in place before the edit.
$ws.Protect()
